$d = $word.ActiveDocument

# Find the paragraph that ends the "Implementando o carrossel" bullet block:
# "Pegamos um dos códigos de exemplo do site e colocamos no nosso index."
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Pegamos um dos códigos de exemplo do site e colocamos no nosso index.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchorPara = $findRange.Paragraphs(1)
$anchorRange = $anchorPara.Range

# --- New paragraph 1 (list level 1): " Quando o resultado não é o que esperamos:" ---
$anchorRange.InsertParagraphAfter()
$p1Range = $d.Paragraphs($d.Paragraphs.Count).Range
$p1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="240" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Quando o resultado não é o que esperamos</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p>'
$p1Range.InsertXML($p1Xml)

# --- New paragraph 2 (list level 2): "Quando isso acontece precisamos olhar o código e ver onde está o problema sempre, isso porque essa também é parte do nosso trabalho." ---
$p2AnchorRange = $d.Paragraphs($d.Paragraphs.Count).Range
$p2AnchorRange.InsertParagraphAfter()
$p2Range = $d.Paragraphs($d.Paragraphs.Count).Range
$p2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="240" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Quando isso acontece precisamos olhar o código e ver onde está o problema sempre, isso porque essa também é parte do nosso trabalho.</w:t></w:r></w:p>'
$p2Range.InsertXML($p2Xml)

Write-Output ("Done. Paragraph count now: " + $d.Paragraphs.Count)
